$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "correo" column header in D1, copying the header style from C1 (bold,
# centered, bordered) so the new column matches the existing header formatting.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("D1").Value = "correo"

# Add the corresponding data value in D2 (no special style, same as other data cells).
$ws.Range("D2").Value = "carlos.velez@est.iudigital.edu.co"
